$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------
# In the "Draft 2" section, the two paragraphs:
#   "This system is implemented entirely as a front-end web application
#    and does not currently include any backend or database
#    integration."                              (red highlight)
#   "It interfaces with the following software libraries and
#    platforms:"                                (yellow highlight)
# are replaced by a single paragraph:
#   "The interface utilizes the following front-end libraries and
#    frameworks:"                               (green highlight)
#
# Find the first paragraph by its text and delete its range (this
# merges it away, leaving the second paragraph intact as the sole
# survivor), then restyle/retext the survivor.

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "This system is implemented entirely*") {
        $target = $p
        break
    }
}
$target.Range.Delete()

$target2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "It interfaces with the following software libraries and platforms:*") {
        $target2 = $p
        break
    }
}
$target2.Range.Text = "The interface utilizes the following front-end libraries and frameworks:"
$target2.Range.HighlightColorIndex = 4

# --- Change 2 --------------------------------------------------------
# Remove the whole paragraph describing GitHub Pages deployment in the
# "Draft 2" section:
#   "The application is deployed using GitHub Pages, which delivers the
#    site as static content. No API endpoints or external services are
#    used in the current version."

$target3 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "The application is deployed using GitHub Pages*") {
        $target3 = $p
        break
    }
}
$target3.Range.Delete()
